# Updated plan for Controller.
#
# The "Control" subsystem section gains a dedicated DBW-Node sub-group:
#   - B25's subsystem label changes from "Control" to "Control – DBW Node"
#   - two new task rows are appended (26: Twist controller, 27: DBW Node)
#     under that same "Control – DBW Node" label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25: subsystem label changes from "Control" to "Control – DBW Node"
$ws.Range("B25").Value = "Control – DBW Node"

# New row 26: Twist controller task
$ws.Range("A26").Value = "Twist controller"
$ws.Range("B26").Value = "Control – DBW Node"
$ws.Range("C26").Value = "This is the main thing which needs to be coded"

# New row 27: DBW Node task
$ws.Range("A27").Value = "DBW Node"
$ws.Range("B27").Value = "Control – DBW Node"
$ws.Range("C27").Value = "Invoke Twist Controller and call publish"

# The "Details" column wraps text for every task row (same look as C12,
# C17, C20, C22 ... C25) - apply the same treatment to the two new rows.
$ws.Range("C26:C27").WrapText = $true

# Match the row height used by the other short task rows (row 25 / row 17
# both use 23.95pt for a 2-line wrapped cell).
$ws.Rows.Item(26).RowHeight = 23.95
$ws.Rows.Item(27).RowHeight = 23.95

# Move the selection onto the newly added detail cells.
$ws.Range("C26:C27").Select()
